$wb = $excel.ActiveWorkbook

$seniorOne = $wb.Worksheets.Item("Senior One")
$seniorOneMarks = @{
    3  = 20
    4  = 60
    5  = 20
    6  = 80
    7  = 80
    8  = 50
    9  = 80
    11 = 40
    12 = 50
    13 = 20
    14 = 30
    16 = 60
    17 = 20
    18 = 50
    22 = 30
    24 = 30
    25 = 80
    27 = 65
    28 = 50
    31 = 60
}
foreach ($row in $seniorOneMarks.Keys) {
    $seniorOne.Range("D$row").Value = $seniorOneMarks[$row]
}

$seniorThree = $wb.Worksheets.Item("Senior Three")
$seniorThreeMarks = @{
    2  = 40
    3  = 80
    4  = 30
    5  = 40
    6  = 70
    7  = 40
    8  = 30
    9  = 30
    12 = 70
    13 = 30
    14 = 40
    15 = 30
    16 = 40
    17 = 30
    18 = 40
    19 = 40
    20 = 80
    21 = 30
    22 = 30
    24 = 70
    26 = 80
    27 = 30
    28 = 80
    29 = 40
    31 = 40
    32 = 70
    33 = 70
    34 = 40
    35 = 80
}
foreach ($row in $seniorThreeMarks.Keys) {
    $seniorThree.Range("D$row").Value = $seniorThreeMarks[$row]
}
